$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the group-member sequence numbers in column A to the new
# 6-digit student-group codes (170300 + old sequence number).
$ws.Range("A3").Value = 170301
$ws.Range("A5").Value = 170302
$ws.Range("A7").Value = 170303
$ws.Range("A10").Value = 170304
$ws.Range("A13").Value = 170305
$ws.Range("A14").Value = 170306
$ws.Range("A15").Value = 170307
$ws.Range("A16").Value = 170308
$ws.Range("A17").Value = 170309
$ws.Range("A18").Value = 170310
$ws.Range("A19").Value = 170311
$ws.Range("A20").Value = 170312
$ws.Range("A21").Value = 170313
$ws.Range("A22").Value = 170314
$ws.Range("A23").Value = 170315
$ws.Range("A24").Value = 170316
$ws.Range("A25").Value = 170317
$ws.Range("A26").Value = 170318
$ws.Range("A27").Value = 170319
$ws.Range("A28").Value = 170320
$ws.Range("A29").Value = 170321

# Move the active selection to A29 (was C31).
$ws.Range("A29").Select()
